$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order (row 1)
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "bedrooms_2"
$ws.Range("E1").Value = "living_rooms_1"
$ws.Range("F1").Value = "living_rooms_2"

# Updated one-hot data rows 2-7
$data = @(
    @(0,0,0,1,0,0),
    @(0,1,0,0,0,0),
    @(0,0,0,0,0,1),
    @(0,0,1,0,0,0),
    @(1,0,0,0,0,0),
    @(0,0,0,0,1,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}

$wb.Save()
